$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Info")
$ws2 = $wb.Worksheets.Item("test_configs")

# Clear the "execute.test" ("x") marks in column A for rows 3-7 of test_configs.
# (Row 2 keeps its "x" - only A3:A7 are being cleared)
$ws2.Range("A3").ClearContents()
$ws2.Range("A4").ClearContents()
$ws2.Range("A5").ClearContents()
$ws2.Range("A6").ClearContents()

$a7 = $ws2.Range("A7")
$a7.ClearContents()
$a7.Borders.LineStyle = -4142
$a7.Font.Name = "Calibri"
$a7.Font.Size = 11
$a7.Interior.Pattern = -4142
$a7.NumberFormat = "GENERAL"
$a7.VerticalAlignment = -4107
$a7.HorizontalAlignment = 1

# Adjust row heights to reflect the new (shorter) content-driven heights
$ws2.Rows.Item(2).RowHeight = 40.95
$ws2.Rows.Item(3).RowHeight = 54.2
$ws2.Rows.Item(4).RowHeight = 54.2
$ws2.Rows.Item(5).RowHeight = 28.1
$ws2.Rows.Item(6).RowHeight = 54.2
$ws2.Rows.Item(7).RowHeight = 28.1

# Update selections: Info sheet selection becomes A3:A7 ...
$ws1.Range("A3:A7").Select()
# ... and test_configs stays/becomes the active sheet again, also selecting A3:A7
$ws2.Activate()
$ws2.Range("A3:A7").Select()

# Shrink the sheet-tab area (tab ratio) of the workbook window
$win = $wb.Windows.Item(1)
$win.TabRatio = 0.293
